$d = $word.ActiveDocument

# Locate the anchor paragraph ("LOQ4057: ...") that must be kept, and the
# last paragraph to be removed (the copyright notice). Everything strictly
# between the end of the anchor paragraph and the end of the copyright
# paragraph (the blank spacer paragraph, the "Ver no Jupiter ..." paragraph,
# and the copyright paragraph itself) is deleted.
$anchorPara = $null
$copyrightPara = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*LOQ4057*") {
        $anchorPara = $p
    }
    if ($t -like "*Powered by Jekyll*") {
        $copyrightPara = $p
    }
}

$rng = $d.Range($anchorPara.Range.End, $copyrightPara.Range.End)
$rng.Delete()
